# Share of Capital and OM Spending by ISIC Code.xlsx
# Split four aggregated ISIC categories on the "SoCaOMSbRIC" sheet into
# their finer-grained sub-codes:
#   ISIC 05T06  -> ISIC 05 | ISIC 06
#   ISIC 23     -> ISIC 231 | ISIC 239
#   ISIC 24     -> ISIC 241 | ISIC 242
#   ISIC 35T39  -> ISIC 351 | ISIC 352T353 | ISIC 36T39
#
# Each split keeps the original column's data in the first sub-column and
# adds new columns (value 0, same header formatting) for the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SoCaOMSbRIC")

# --- Split "ISIC 05T06" (column C) into "ISIC 05" and "ISIC 06" ---
$ws.Columns("C").Copy()
$ws.Columns("D").Insert()
$ws.Range("C1").Value = "ISIC 05"
$ws.Range("D1").Value = "ISIC 06"
$ws.Range("D2").Value = 0

# --- Split "ISIC 23" (now column O) into "ISIC 231" and "ISIC 239" ---
$ws.Columns("O").Copy()
$ws.Columns("P").Insert()
$ws.Range("O1").Value = "ISIC 231"
$ws.Range("P1").Value = "ISIC 239"
$ws.Range("P2").Value = 0

# --- Split "ISIC 24" (now column Q) into "ISIC 241" and "ISIC 242" ---
$ws.Columns("Q").Copy()
$ws.Columns("R").Insert()
$ws.Range("Q1").Value = "ISIC 241"
$ws.Range("R1").Value = "ISIC 242"
$ws.Range("R2").Value = 0

# --- Split "ISIC 35T39" (now column Z) into "ISIC 351", "ISIC 352T353", "ISIC 36T39" ---
$ws.Columns("Z").Copy()
$ws.Columns("AA").Insert()
$ws.Columns("Z").Copy()
$ws.Columns("AB").Insert()
$ws.Range("Z1").Value = "ISIC 351"
$ws.Range("AA1").Value = "ISIC 352T353"
$ws.Range("AB1").Value = "ISIC 36T39"
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
